$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.573.92'
$ws.Range('D2').NumberFormat = "General"
$ws.Range('D2').Style = "Normal"

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('E2').NumberFormat = "General"
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.798.20'
$ws.Range('D3').NumberFormat = "General"
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E3').NumberFormat = "General"
$ws.Range('E3').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.35'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('E5').NumberFormat = "General"
$ws.Range('E5').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E7').NumberFormat = "General"
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.87'
$ws.Range('D8').NumberFormat = "General"
$ws.Range('D8').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +3.67%  '
$ws.Range('E8').NumberFormat = "General"
$ws.Range('E8').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('E9').NumberFormat = "General"
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0694'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('D10').Style = "Normal"

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('E10').NumberFormat = "General"
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0949'
$ws.Range('D11').NumberFormat = "General"
$ws.Range('D11').Style = "Normal"

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E11').NumberFormat = "General"
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.058.97'
$ws.Range('D12').NumberFormat = "General"
$ws.Range('D12').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('E12').NumberFormat = "General"
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.12'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('D13').Style = "Normal"

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('E13').NumberFormat = "General"
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.801.35'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('D14').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.44%  '
$ws.Range('E14').NumberFormat = "General"
$ws.Range('E14').Style = "Normal"

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('E15').NumberFormat = "General"
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '34.574.93'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('D16').Style = "Normal"

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('E16').NumberFormat = "General"
$ws.Range('E16').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +3.15%  '
$ws.Range('E17').NumberFormat = "General"
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.87'
$ws.Range('D18').NumberFormat = "General"
$ws.Range('D18').Style = "Normal"

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('E19').NumberFormat = "General"
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '246.55'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('D20').Style = "Normal"

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E20').NumberFormat = "General"
$ws.Range('E20').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +3.67%  '
$ws.Range('E21').NumberFormat = "General"
$ws.Range('E21').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.16'
$ws.Range('D23').NumberFormat = "General"
$ws.Range('D23').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E23').NumberFormat = "General"
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '173.48'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +6.76%  '
$ws.Range('E24').NumberFormat = "General"
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.07'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D25').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('E25').NumberFormat = "General"
$ws.Range('E25').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('E26').NumberFormat = "General"
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.67'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E27').NumberFormat = "General"
$ws.Range('E27').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('E28').NumberFormat = "General"
$ws.Range('E28').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E29').NumberFormat = "General"
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.03'
$ws.Range('D30').NumberFormat = "General"
$ws.Range('D30').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +8.39%  '
$ws.Range('E30').NumberFormat = "General"
$ws.Range('E30').Style = "Normal"

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('E31').NumberFormat = "General"
$ws.Range('E31').Style = "Normal"

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('E32').NumberFormat = "General"
$ws.Range('E32').Style = "Normal"

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('E33').NumberFormat = "General"
$ws.Range('E33').Style = "Normal"

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('E34').NumberFormat = "General"
$ws.Range('E34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.431.59'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('D35').Style = "Normal"

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('E35').NumberFormat = "General"
$ws.Range('E35').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +6.66%  '
$ws.Range('E36').NumberFormat = "General"
$ws.Range('E36').Style = "Normal"

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.676'
$ws.Range('D37').NumberFormat = "General"
$ws.Range('D37').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +2.28%  '
$ws.Range('E37').NumberFormat = "General"
$ws.Range('E37').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('E38').NumberFormat = "General"
$ws.Range('E38').Style = "Normal"

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('E39').NumberFormat = "General"
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '84.73'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('D40').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +5.64%  '
$ws.Range('E40').NumberFormat = "General"
$ws.Range('E40').Style = "Normal"

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('E41').NumberFormat = "General"
$ws.Range('E41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.40'
$ws.Range('D42').NumberFormat = "General"
$ws.Range('D42').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('E42').NumberFormat = "General"
$ws.Range('E42').Style = "Normal"

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('E43').NumberFormat = "General"
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.88'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('D44').Style = "Normal"

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +3.24%  '
$ws.Range('E44').NumberFormat = "General"
$ws.Range('E44').Style = "Normal"

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('E45').NumberFormat = "General"
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('E46').NumberFormat = "General"
$ws.Range('E46').Style = "Normal"

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('E47').NumberFormat = "General"
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.959.28'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('E48').NumberFormat = "General"
$ws.Range('E48').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '105.20'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('D49').Style = "Normal"

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('E49').NumberFormat = "General"
$ws.Range('E49').Style = "Normal"

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('E50').NumberFormat = "General"
$ws.Range('E50').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -5.04%  '
$ws.Range('E51').NumberFormat = "General"
$ws.Range('E51').Style = "Normal"

